$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("G2").Value = 92.39013433333332
$ws.Range("H2").Value = 277.170403
$ws.Range("I2").Value = 0.2557671247968706
$ws.Range("J2").Value = 0.2557671247968705
$ws.Range("K2").Value = 3
$ws.Range("L2").Value = 1
$ws.Range("M2").Value = 0.4910443333333334
$ws.Range("N2").Value = 1.473133
$ws.Range("O2").Value = 0.7844104380534107
$ws.Range("P2").Value = 0.7844104380534107
$ws.Range("Q2").Value = 45.36765192028877
$ws.Range("R2").Value = 408.308867282599
$ws.Range("S2").Value = 0.2006264024015746
$ws.Range("T2").Value = 0.2006264024015745
$ws.Range("G3").Value = 92.39013433333332
$ws.Range("H3").Value = 277.170403
$ws.Range("I3").Value = 0.2557671247968706
$ws.Range("J3").Value = 0.2557671247968705
$ws.Range("K3").Value = 2
$ws.Range("L3").Value = 0.6666666666666666
$ws.Range("M3").Value = 0.13496
$ws.Range("N3").Value = 0.40488
$ws.Range("O3").Value = 0.2155895619465893
$ws.Range("P3").Value = 0.2155895619465893
$ws.Range("Q3").Value = 12.46897252962667
$ws.Range("R3").Value = 112.22075276664
$ws.Range("S3").Value = 0.05514072239529596
$ws.Range("T3").Value = 0.05514072239529594
$ws.Range("I4").Value = 0.7056897640051698
$ws.Range("J4").Value = 0.7056897640051695
$ws.Range("K4").Value = 3
$ws.Range("L4").Value = 1
$ws.Range("M4").Value = 0.4910443333333334
$ws.Range("N4").Value = 1.473133
$ws.Range("O4").Value = 0.7844104380534107
$ws.Range("P4").Value = 0.7844104380534107
$ws.Range("Q4").Value = 125.1743655582158
$ws.Range("R4").Value = 1126.569290023942
$ws.Range("S4").Value = 0.5535504169131032
$ws.Range("T4").Value = 0.553550416913103
$ws.Range("I5").Value = 0.7056897640051698
$ws.Range("J5").Value = 0.7056897640051695
$ws.Range("K5").Value = 2
$ws.Range("L5").Value = 0.6666666666666666
$ws.Range("M5").Value = 0.13496
$ws.Range("N5").Value = 0.40488
$ws.Range("O5").Value = 0.2155895619465893
$ws.Range("P5").Value = 0.2155895619465893
$ws.Range("Q5").Value = 34.40327324634666
$ws.Range("R5").Value = 309.62945921712
$ws.Range("S5").Value = 0.1521393470920665
$ws.Range("T5").Value = 0.1521393470920665
$ws.Range("E6").Value = 3
$ws.Range("F6").Value = 1
$ws.Range("G6").Value = 0.3441203333333334
$ws.Range("H6").Value = 1.032361
$ws.Range("I6").Value = 0.0009526414143230948
$ws.Range("J6").Value = 0.0009526414143230943
$ws.Range("K6").Value = 3
$ws.Range("L6").Value = 1
$ws.Range("M6").Value = 0.4910443333333334
$ws.Range("N6").Value = 1.473133
$ws.Range("O6").Value = 0.7844104380534107
$ws.Range("P6").Value = 0.7844104380534107
$ws.Range("Q6").Value = 0.1689783396681111
$ws.Range("R6").Value = 1.520805057013
$ws.Range("S6").Value = 0.0007472618691169995
$ws.Range("T6").Value = 0.0007472618691169991
$ws.Range("E7").Value = 3
$ws.Range("F7").Value = 1
$ws.Range("G7").Value = 0.3441203333333334
$ws.Range("H7").Value = 1.032361
$ws.Range("I7").Value = 0.0009526414143230948
$ws.Range("J7").Value = 0.0009526414143230943
$ws.Range("K7").Value = 2
$ws.Range("L7").Value = 0.6666666666666666
$ws.Range("M7").Value = 0.13496
$ws.Range("N7").Value = 0.40488
$ws.Range("O7").Value = 0.2155895619465893
$ws.Range("P7").Value = 0.2155895619465893
$ws.Range("Q7").Value = 0.04644248018666667
$ws.Range("R7").Value = 0.41798232168
$ws.Range("S7").Value = 0.0002053795452060952
$ws.Range("T7").Value = 0.0002053795452060952
$ws.Range("G8").Value = 13.48925766666667
$ws.Range("H8").Value = 40.467773
$ws.Range("I8").Value = 0.0373428253345738
$ws.Range("J8").Value = 0.03734282533457379
$ws.Range("K8").Value = 3
$ws.Range("L8").Value = 1
$ws.Range("M8").Value = 0.4910443333333334
$ws.Range("N8").Value = 1.473133
$ws.Range("O8").Value = 0.7844104380534107
$ws.Range("P8").Value = 0.7844104380534107
$ws.Range("Q8").Value = 6.62382353808989
$ws.Range("R8").Value = 59.61441184280901
$ws.Range("S8").Value = 0.02929210197884504
$ws.Range("T8").Value = 0.02929210197884503
$ws.Range("G9").Value = 13.48925766666667
$ws.Range("H9").Value = 40.467773
$ws.Range("I9").Value = 0.0373428253345738
$ws.Range("J9").Value = 0.03734282533457379
$ws.Range("K9").Value = 2
$ws.Range("L9").Value = 0.6666666666666666
$ws.Range("M9").Value = 0.13496
$ws.Range("N9").Value = 0.40488
$ws.Range("O9").Value = 0.2155895619465893
$ws.Range("P9").Value = 0.2155895619465893
$ws.Range("Q9").Value = 1.820510214693333
$ws.Range("R9").Value = 16.38459193224
$ws.Range("S9").Value = 0.008050723355728762
$ws.Range("T9").Value = 0.00805072335572876
$ws.Range("G10").Value = 0.08945599999999999
$ws.Range("H10").Value = 0.268368
$ws.Range("I10").Value = 0.0002476444490629346
$ws.Range("J10").Value = 0.0002476444490629346
$ws.Range("K10").Value = 3
$ws.Range("L10").Value = 1
$ws.Range("M10").Value = 0.4910443333333334
$ws.Range("N10").Value = 1.473133
$ws.Range("O10").Value = 0.7844104380534107
$ws.Range("P10").Value = 0.7844104380534107
$ws.Range("Q10").Value = 0.04392686188266667
$ws.Range("R10").Value = 0.395341756944
$ws.Range("S10").Value = 0.0001942548907709521
$ws.Range("T10").Value = 0.000194254890770952
$ws.Range("G11").Value = 0.08945599999999999
$ws.Range("H11").Value = 0.268368
$ws.Range("I11").Value = 0.0002476444490629346
$ws.Range("J11").Value = 0.0002476444490629346
$ws.Range("K11").Value = 2
$ws.Range("L11").Value = 0.6666666666666666
$ws.Range("M11").Value = 0.13496
$ws.Range("N11").Value = 0.40488
$ws.Range("O11").Value = 0.2155895619465893
$ws.Range("P11").Value = 0.2155895619465893
$ws.Range("Q11").Value = 0.01207298176
$ws.Range("R11").Value = 0.10865683584
$ws.Range("S11").Value = 0.00005338955829198251
$ws.Range("T11").Value = 0.00005338955829198251
